$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "ser" blog identifiers up by one: 130 -> 131, 131 -> 132, 132 -> 133
$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 131"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 132"
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 133"
